# windata10m.xlsx — add a "Power delivered" half-hourly share column (O),
# point the existing daily-average formula at the new column, and keep the
# old N-based daily average around in a new R3 cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("windata10m")

# --- New column O: O{r} = N{r} * 1/6 for every data row (2..145) ---------
# Row 2 is written on its own (matches the source file, which has it as a
# standalone formula rather than part of a shared-formula block), while the
# rest are written in the same three blocks that column N already uses so
# the shared-formula grouping lines up with N3:N66 / N67:N130 / N131:N145.
$ws.Range("O2").Formula = "=N2*1/6"
$ws.Range("O3:O66").Formula = "=N3*1/6"
$ws.Range("O67:O130").Formula = "=N67*1/6"
$ws.Range("O131:O145").Formula = "=N131*1/6"

# --- R3 keeps the original daily-average-from-N formula -------------------
$ws.Range("R3").Formula = "=SUM(N2:N145)/145*24"

# --- P3 now sums the new O (per-interval share) column instead -----------
$ws.Range("P3").Formula = "=SUM(O2:O145)"

# --- Column A: widen to fit the timestamps --------------------------------
$ws.Columns.Item(1).ColumnWidth = 16.5

# --- Selection moves from P4 to P3 ----------------------------------------
$ws.Range("P3").Select()
